$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for retired point-release .NET Framework versions, starting
# from the bottom so earlier row numbers stay valid.
$ws.Rows.Item(16).Delete()   # .NET Framework 4.7.2
$ws.Rows.Item(15).Delete()   # .NET Framework 4.7.1
$ws.Rows.Item(13).Delete()   # .NET Framework 4.6.2
$ws.Rows.Item(12).Delete()   # .NET Framework 4.6.1
$ws.Rows.Item(10).Delete()   # .NET Framework 4.5.2
$ws.Rows.Item(9).Delete()    # .NET Framework 4.5.1

# Give the two new rows the same formatting as the existing data rows
$ws.Range("A19:B19").Copy()
$ws.Range("A20:B21").PasteSpecial(-4122)  # xlPasteFormats

# Append new rows for .NET 6 and .NET 7 (write labels first, then dates,
# so new shared strings are interned in the same order as the target file)
$ws.Cells.Item(20, 1).Value = ".NET 6"
$ws.Cells.Item(21, 1).Value = ".NET 7"
$ws.Cells.Item(20, 2).Value = "November 8, 2021"
$ws.Cells.Item(21, 2).Value = "November 8, 2022"

# Update the selection to match the target state
$ws.Range("A2").Select()
